# "Add files via upload" - adds 4 new rows to the "r0" sheet documenting
# Hanfei's flow-rate optimization experiment (script variants run at
# 2, 5, 10 and 25 mL/hr sample flow rates).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("r0")

$name = "Hanfei's flow rate opttimization"

# Common values shared by every new row, keyed by column letter.
$common = @{
    "C" = "20 mL"          # Waste Syringe
    "D" = "5 mL"            # Lysate Syringe
    "E" = "1 hour"          # F-127 Incubation Time
    "F" = "2.5 mL"          # Sample Volume
    "H" = "15 mL/hr"        # PBS wash flow rate (post-sample)
    "I" = "200-800-1000"    # PBS Wash Structure
    "J" = "2 mins"          # QIAzol Incubation Time
    "K" = "N"               # PBS wash post-QIAzol?
}

# Row number -> Sample Flow Rate (column G) for each new script variant.
$rows = [ordered]@{
    21 = "2 mL/hr"
    22 = "5 mL/hr"
    23 = "10 mL/hr"
    24 = "25 mL/hr"
}

foreach ($r in $rows.Keys) {
    $ws.Range("B$r").Value = $name
    foreach ($col in $common.Keys) {
        $ws.Range("$col$r").Value = $common[$col]
    }
    $ws.Range("G$r").Value = $rows[$r]
}

# Row 24 carries an explicit (slightly taller) custom row height.
$ws.Rows.Item(24).RowHeight = 17

$ws.Activate()
$ws.Range("B24").Select()
